# BUG: Don't extract header names if none specified (pandas#23703)
# Adds a new worksheet "index_col_none" used as test fixture data for
# pandas' read_excel with a MultiIndex column header and no index_col.
#
# Layout (4 cols x 4 rows):
#   Row1:  A    A    B    B       (bold, centered, no border)
#   Row2:  key  val  key  val     (bold, centered, no border)
#   Row3:  1    2    3    4       (centered, no border)
#   Row4:  1    2    3    4       (centered, no border)

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "index_col_none"

# Header rows (two-level MultiIndex column header).
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "B"

$ws.Range("A2").Value = "key"
$ws.Range("B2").Value = "val"
$ws.Range("C2").Value = "key"
$ws.Range("D2").Value = "val"

# Data rows.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4

# Formatting: data rows centered first, then header rows centered + bold,
# so the header's bold style reuses the existing bold/centered/no-border
# style already present in the workbook instead of spawning an extra xf.
$dataRange = $ws.Range("A3:D4")
$dataRange.HorizontalAlignment = -4108  # xlCenter

$headerRange = $ws.Range("A1:D2")
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.Font.Bold = $true
